$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column O (DÉLKA_PRACOVNÍHO_POMĚRU) first so column M's index is not affected,
# then delete column M (ZAŘAZENO).
$ws.Range("O1").EntireColumn.Delete()
$ws.Range("M1").EntireColumn.Delete()

# Reflect the resulting selection/scroll state seen after deleting the columns.
$null = $ws.Range("M1:M1048576").Select()
$excel.ActiveWindow.ScrollColumn = 12
